$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: id / value
$ws.Range("G21").Value = "XpOLJXtn"
$ws.Range("F21").Value = "id"
$ws.Range("F21").Font.Bold = $true

# Row 22: duration (hours)
$ws.Range("F22").Value = "duration (hours)"
$ws.Range("F22").Font.Bold = $true
$ws.Range("H22").Value = 0.97986111111111107
$ws.Range("G22").Formula = "=H22*24"
$ws.Range("H22").NumberFormat = "h:mm:ss"

# Row 23: trials
$ws.Range("F23").Value = "trials"
$ws.Range("F23").Font.Bold = $true
$ws.Range("G23").Value = 658

# Row 24: best result
$ws.Range("F24").Value = "best result"
$ws.Range("F24").Font.Bold = $true
$ws.Range("G24").Value = 0.96194299999999999

[void]$ws.Range("F21:G24").Select()
